$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Target cluster "ECs" instead of "FAPs", plus new metric values) ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.424501
$ws.Range("H2").Value = 10.273503
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.033285
$ws.Range("N2").Value = 0.099855
$ws.Range("O2").Value = 0.007684499559038781
$ws.Range("P2").Value = 0.007684499559038781
$ws.Range("Q2").Value = 0.113984515785
$ws.Range("R2").Value = 1.025860642065
$ws.Range("S2").Value = 0.007684499559038781
$ws.Range("T2").Value = 0.007684499559038781

# --- Update existing row 3 to hold the "FAPs" target-cluster data with new metric values ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf9"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.424501
$ws.Range("H3").Value = 10.273503
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.140873
$ws.Range("N3").Value = 12.422619
$ws.Range("O3").Value = 0.9560023056192156
$ws.Range("P3").Value = 0.9560023056192157
$ws.Range("Q3").Value = 14.180423729373
$ws.Range("R3").Value = 127.623813564357
$ws.Range("S3").Value = 0.9560023056192156
$ws.Range("T3").Value = 0.9560023056192157

# --- Add new row 4 for the "sCs" target-cluster data ---
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf9"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.424501
$ws.Range("H4").Value = 10.273503
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1572886666666667
$ws.Range("N4").Value = 0.471866
$ws.Range("O4").Value = 0.03631319482174546
$ws.Range("P4").Value = 0.03631319482174546
$ws.Range("Q4").Value = 0.5386351962886666
$ws.Range("R4").Value = 4.847716766598
$ws.Range("S4").Value = 0.03631319482174546
$ws.Range("T4").Value = 0.03631319482174546
